# Weekly fruit/vegetable data update:
# A new observation is inserted as row 602 (shifting the existing rows
# 602-647 down to 603-648), growing the used range from A1:R647 to A1:R648.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 602, pushing everything
# below it (through row 647) down by one row.
$ws.Rows(602).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A602").Value = 6
$ws.Range("B602").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C602").Value = "Metropolitana"
$ws.Range("D602").Value = 45132
$ws.Range("E602").Value = 13
$ws.Range("F602").Value = 100112043
$ws.Range("G602").Value = "Pepino ensalada"
$ws.Range("H602").Value = "Sin especificar"
$ws.Range("I602").Value = "Primera"
$ws.Range("J602").Value = 230
$ws.Range("K602").Value = 8000
$ws.Range("L602").Value = 8000
$ws.Range("M602").Value = 8000
$ws.Range("N602").Value = "$/caja 60 unidades"
$ws.Range("O602").Value = "Región de Arica y Parinacota"
$ws.Range("P602").Value = 133
$ws.Range("Q602").Value = 60
$ws.Range("R602").Value = "Hortaliza"
